# PropertyTypes.xlsx: "StringNull of" row used to show the plain text `null`
# in column D. The fix now emits the quoted string `"null"` instead, so the
# sheet has to reflect a quoted-string sample value for that row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 / column D holds the sample value for "StringNull of". Replace the
# bare word null with the literal text "null" (including the quote marks).
$ws.Range("D6").Value = """null"""

# Leave the selection where it ended up after making the edit.
$ws.Range("D7").Select()
